# Poster A1, Bitmessage for Android - update the Bitmessage address sample
# shown on the slide and let the text box settle into its new (slightly
# taller) auto-fit height.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Textfeld 5" shape that holds the "What is it about?" copy,
# including the sample Bitmessage address in its last paragraph.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Textfeld 5") {
        $shp = $candidate
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(5)
}

$tr = $shp.TextFrame.TextRange

$oldAddress = "BM-2cWs84ik1Fj7jdJKrn3vDecxQbH9R4VS9r"
$newSuffix = "-2cUau5uxBYCK2Z2TVwUZnnNfYW5yyutekC"

$fullText = $tr.Text
$addressStart = $fullText.IndexOf($oldAddress)

# Keep the leading "BM" untouched and only retype the remainder of the
# address (mirrors how the address was edited by hand in PowerPoint).
$suffixRange = $tr.Characters($addressStart + 3, $oldAddress.Length - 2)
$suffixRange.Text = $newSuffix

# The text box auto-fits its height to the paragraph content
# (<a:spAutoFit/>); resolve the new laid-out extent in points so the
# stored EMU height matches PowerPoint's own re-layout after the edit.
$targetHeightEmu = 10495182
$shp.Height = ($targetHeightEmu / 12700.0) + 0.00005
